# Update "想去人数" (F column) figures for the refreshed scrape (456a3b4).
$wb = $excel.ActiveWorkbook

# Sheet "展览" (index 1)
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F2").Value = 768
$ws1.Range("F6").Value = 118
$ws1.Range("F8").Value = 123
$ws1.Range("F9").Value = 324
$ws1.Range("F10").Value = 432
$ws1.Range("F11").Value = 498
$ws1.Range("F13").Value = 11415
$ws1.Range("F14").Value = 5375

# Sheet "演出" (index 2)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F2").Value = 96

# Sheet "全部类型" (index 4)
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F2").Value = 768
$ws4.Range("F5").Value = 96
$ws4.Range("F8").Value = 118
$ws4.Range("F10").Value = 123
$ws4.Range("F11").Value = 324
$ws4.Range("F12").Value = 432
$ws4.Range("F13").Value = 498
$ws4.Range("F15").Value = 11415
$ws4.Range("F17").Value = 5375
